$wb = $excel.ActiveWorkbook

# Update the "Status" text wherever it currently reads "Ready for handoff" so it
# reflects that the item is now in translation. This shared string is used by
# the Overview sheet as well as each per-language sheet.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# The status columns were sized to fit the old, longer text ("Ready for
# handoff"). Now that the text is shorter ("In Translation") re-fit those
# columns so the report looks right when it is generated for archive.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns("C:C").ColumnWidth = 12.5
